$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 8990.134
$ws.Range("I137").Value = 10487.75
$ws.Range("J137").Value = 2999.6667
$ws.Range("K137").Value = 31463.25
$ws.Range("L137").Value = 8999.000100000001
$ws.Range("M137").Value = -28913.25
$ws.Range("N137").Value = -14099.0001
$ws.Range("H138").Value = 2280.196
$ws.Range("J138").Value = 2763.8857
$ws.Range("L138").Value = 8291.6571
$ws.Range("N138").Value = -18571.6571

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 145.83333
$ws.Range("I5").Value = 95
$ws.Range("K5").Value = 95
$ws.Range("M5").Value = 17
$ws.Range("H17").Value = 3000
$ws.Range("I17").Value = 3000
$ws.Range("K17").Value = 3000
$ws.Range("M17").Value = -2827
$ws.Range("H61").Value = 4721.909
$ws.Range("I61").Value = 5049.722
$ws.Range("K61").Value = 5049.722
$ws.Range("M61").Value = -4837.722
$ws.Range("H63").Value = 4353.4287
$ws.Range("J63").Value = 8360
$ws.Range("L63").Value = 8360
$ws.Range("N63").Value = -9732
$ws.Range("H66").Value = 4353.4287
$ws.Range("J66").Value = 8360
$ws.Range("L66").Value = 41800
$ws.Range("N66").Value = -48664
$ws.Range("H74").Value = 2307.1
$ws.Range("I74").Value = 2376.2104
$ws.Range("J74").Value = 994
$ws.Range("K74").Value = 2376.2104
$ws.Range("L74").Value = 994
$ws.Range("M74").Value = -1502.2104
$ws.Range("N74").Value = -2742
$ws.Range("H77").Value = 2307.1
$ws.Range("I77").Value = 2376.2104
$ws.Range("J77").Value = 994
$ws.Range("K77").Value = 11881.052
$ws.Range("L77").Value = 4970
$ws.Range("M77").Value = -7513.052
$ws.Range("N77").Value = -13706
$ws.Range("H88").Value = 1385.9375
$ws.Range("I88").Value = 498.5
$ws.Range("J88").Value = 1918.4
$ws.Range("K88").Value = 498.5
$ws.Range("L88").Value = 1918.4
$ws.Range("M88").Value = -92.5
$ws.Range("N88").Value = -2730.4
$ws.Range("H91").Value = 1385.9375
$ws.Range("I91").Value = 498.5
$ws.Range("J91").Value = 1918.4
$ws.Range("K91").Value = 498.5
$ws.Range("L91").Value = 1918.4
$ws.Range("M91").Value = 905.5
$ws.Range("N91").Value = -4726.4
$ws.Range("H132").Value = 2756.7273
$ws.Range("I132").Value = 1916.125
$ws.Range("K132").Value = 5748.375
$ws.Range("M132").Value = -3218.375
$ws.Range("H136").Value = 4721.909
$ws.Range("I136").Value = 5049.722
$ws.Range("K136").Value = 15149.166
$ws.Range("M136").Value = -12599.166

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 145.83333
$ws.Range("I4").Value = 95
$ws.Range("K4").Value = 95
$ws.Range("M4").Value = 20
$ws.Range("H22").Value = 1372.2858
$ws.Range("I22").Value = 1180.9
$ws.Range("K22").Value = 1180.9
$ws.Range("M22").Value = -1007.9
$ws.Range("H86").Value = 2757.9167
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 2757.9167
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H135").Value = 64999.4
$ws.Range("J135").Value = 64999.4
$ws.Range("L135").Value = 64999.4
$ws.Range("N135").Value = -75139.39999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5513.75
$ws.Range("I31").Value = 1935.25
$ws.Range("J31").Value = 16249.25
$ws.Range("K31").Value = 1935.25
$ws.Range("L31").Value = 16249.25
$ws.Range("M31").Value = -1640.25
$ws.Range("N31").Value = -16839.25
$ws.Range("H34").Value = 5513.75
$ws.Range("I34").Value = 1935.25
$ws.Range("J34").Value = 16249.25
$ws.Range("K34").Value = 1935.25
$ws.Range("L34").Value = 16249.25
$ws.Range("M34").Value = -1733.25
$ws.Range("N34").Value = -16653.25
$ws.Range("H99").Value = 3725
$ws.Range("I99").Value = 3087.5
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 3087.5
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = -1589.5
$ws.Range("N99").Value = -7996
$ws.Range("H126").Value = 3725
$ws.Range("I126").Value = 3087.5
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 9262.5
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -6792.5
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 2915.111
$ws.Range("I132").Value = 2942.3845
$ws.Range("J132").Value = 2844.2
$ws.Range("K132").Value = 8827.1535
$ws.Range("L132").Value = 8532.599999999999
$ws.Range("M132").Value = -6297.1535
$ws.Range("N132").Value = -13592.6

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1252.5834
$ws.Range("I5").Value = 440.14285
$ws.Range("K5").Value = 1320.42855
$ws.Range("M5").Value = -1208.42855
$ws.Range("H68").Value = 4610
$ws.Range("J68").Value = 5061.375
$ws.Range("L68").Value = 15184.125
$ws.Range("N68").Value = -16806.125
$ws.Range("H71").Value = 4610
$ws.Range("J71").Value = 5061.375
$ws.Range("L71").Value = 45552.375
$ws.Range("N71").Value = -53664.375
$ws.Range("H97").Value = 810.2857
$ws.Range("I97").Value = 174.2
$ws.Range("J97").Value = 1163.6666
$ws.Range("K97").Value = 522.5999999999999
$ws.Range("L97").Value = 3490.9998
$ws.Range("M97").Value = -26.59999999999991
$ws.Range("N97").Value = -4482.9998
$ws.Range("H113").Value = 1123.5294
$ws.Range("J113").Value = 1176.6666
$ws.Range("L113").Value = 3529.9998
$ws.Range("N113").Value = -7869.9998
$ws.Range("H120").Value = 120842.71
$ws.Range("I120").Value = 507904
$ws.Range("K120").Value = 1523712
$ws.Range("M120").Value = -1518874
$ws.Range("H132").Value = 2790.818
$ws.Range("I132").Value = 2233.3333
$ws.Range("K132").Value = 20099.9997
$ws.Range("M132").Value = -17569.9997
$ws.Range("H135").Value = 1252.5834
$ws.Range("I135").Value = 440.14285
$ws.Range("K135").Value = 3961.28565
$ws.Range("M135").Value = -1426.28565
$ws.Range("H137").Value = 4867.8
$ws.Range("I137").Value = 1332.4
$ws.Range("J137").Value = 11938.6
$ws.Range("K137").Value = 3997.2
$ws.Range("L137").Value = 35815.8
$ws.Range("M137").Value = 1102.8
$ws.Range("N137").Value = -46015.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 12991.917
$ws.Range("I80").Value = 18129
$ws.Range("K80").Value = 18129
$ws.Range("M80").Value = -17131
$ws.Range("H83").Value = 12991.917
$ws.Range("I83").Value = 18129
$ws.Range("K83").Value = 90645
$ws.Range("M83").Value = -85653
$ws.Range("H132").Value = 2487.6
$ws.Range("I132").Value = 2966.1333
$ws.Range("K132").Value = 8898.3999
$ws.Range("M132").Value = -6368.3999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 4280
$ws.Range("I82").Value = 3011.4285
$ws.Range("J82").Value = 6500
$ws.Range("K82").Value = 3011.4285
$ws.Range("L82").Value = 6500
$ws.Range("M82").Value = -2650.4285
$ws.Range("N82").Value = -7222
$ws.Range("H85").Value = 4280
$ws.Range("I85").Value = 3011.4285
$ws.Range("J85").Value = 6500
$ws.Range("K85").Value = 3011.4285
$ws.Range("L85").Value = 6500
$ws.Range("M85").Value = -1763.4285
$ws.Range("N85").Value = -8996
$ws.Range("H93").Value = 3533.4285
$ws.Range("I93").Value = 1727.1666
$ws.Range("J93").Value = 4888.125
$ws.Range("K93").Value = 1727.1666
$ws.Range("L93").Value = 4888.125
$ws.Range("M93").Value = -479.1666
$ws.Range("N93").Value = -7384.125

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 10076
$ws.Range("J69").Value = 10076
$ws.Range("L69").Value = 10076
$ws.Range("N69").Value = -11574
$ws.Range("H72").Value = 10076
$ws.Range("J72").Value = 10076
$ws.Range("L72").Value = 30228
$ws.Range("N72").Value = -37716
$ws.Range("H81").Value = 2600
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H84").Value = 2600
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H107").Value = 4550.364
$ws.Range("J107").Value = 14744.75
$ws.Range("L107").Value = 44234.25
$ws.Range("N107").Value = -48074.25
$ws.Range("H113").Value = 1814.6666
$ws.Range("I113").Value = 1703.4445
$ws.Range("J113").Value = 2148.3333
$ws.Range("K113").Value = 5110.333500000001
$ws.Range("L113").Value = 6444.999899999999
$ws.Range("M113").Value = -2940.333500000001
$ws.Range("N113").Value = -10784.9999
